{"js": "// The \"Assignment 2: Checklist\" section has four tables (Pass, 2:2, 2:1,\n// First standard). Each row is a grading criterion followed by a Yes/No\n// answer cell. Adding the Statistics class let several previously unmet\n// criteria pass, so their answer cell flips from \"No\" to \"Yes\".\n//\n// (tableIndex, rowIndex) pairs identify the criterion row (0-based) whose\n// answer cell (column 1) changes from \"No\" to \"Yes\":\nconst changes = [\n  [0, 2], // Object instantiation, method calls evident.\n  [0, 3], // Sevens Out game is created.\n  [1, 0], // The rules of the Sevens Out game, as specified, are implemented.\n  [1, 1], // Application repeats or quits the game gracefully according to user choice.\n  [1, 2], // Method calls from 'Main' to methods in other classes\n  [1, 3], // Error handling is evident, some errors are captured, such as erroneous input being made.\n  [1, 4], // Class definitions show encapsulation.\n  [2, 0], // Sevens Out and Three Or More games are implemented.\n  [2, 1], // Inheritance is implemented, showing a class hierarchy\n  [2, 2], // public/private access control in classes\n  [2, 3], // Generic collections (such as List<>) are used.\n  [2, 4], // Exception handling is used\n  [3, 0], // Interfaces and LINQ are used\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst cells = changes.map(([tableIndex, rowIndex]) =>\n  tables.items[tableIndex].getCell(rowIndex, 1)\n);\ncells.forEach((cell) => cell.body.load(\"text\"));\nawait context.sync();\n\ncells.forEach((cell) => {\n  if (cell.body.text.trim() === \"No\") {\n    cell.body.insertText(\"Yes\", Word.InsertLocation.replace);\n  }\n});\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# The \"Assignment 2: Checklist\" section has four tables (Pass, 2:2, 2:1,\n# First standard). Each row is a grading criterion followed by a Yes/No\n# answer cell. Adding the Statistics class let several previously unmet\n# criteria pass, so their answer cell flips from \"No\" to \"Yes\".\n#\n# (tableIndex, rowIndex) are 1-based Word COM Tables/Rows indices; column 2\n# holds the Yes/No answer for the criterion described in column 1.\n$changes = @(\n  @(1, 3),  # Object instantiation, method calls evident.\n  @(1, 4),  # Sevens Out game is created.\n  @(2, 1),  # The rules of the Sevens Out game, as specified, are implemented.\n  @(2, 2),  # Application repeats or quits the game gracefully according to user choice.\n  @(2, 3),  # Method calls from 'Main' to methods in other classes\n  @(2, 4),  # Error handling is evident, some errors are captured, such as erroneous input being made.\n  @(2, 5),  # Class definitions show encapsulation.\n  @(3, 1),  # Sevens Out and Three Or More games are implemented.\n  @(3, 2),  # Inheritance is implemented, showing a class hierarchy\n  @(3, 3),  # public/private access control in classes\n  @(3, 4),  # Generic collections (such as List<>) are used.\n  @(3, 5),  # Exception handling is used\n  @(4, 1)   # Interfaces and LINQ are used\n)\n\nforeach ($change in $changes) {\n  $tableIndex = $change[0]\n  $rowIndex = $change[1]\n  $cell = $d.Tables.Item($tableIndex).Cell($rowIndex, 2)\n  if ($cell.Range.Text -like \"No*\") {\n    $cell.Range.Text = \"Yes\"\n  }\n}\n"}
